# Update TCR, starcat, and various smaller changes
# Correct the relapse sample-id suffixes in column D from lowercase "_rel"
# (and one stray "_Rem" typo) to the consistent "_Rel" used elsewhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D98").Value = "P27_Rel"
$ws.Range("D101").Value = "P28_Rel"
$ws.Range("D104").Value = "P29_Rel"
$ws.Range("D105").Value = "P29_Rel"
$ws.Range("D110").Value = "P30_Rel"
$ws.Range("D115").Value = "P31_Rel"
$ws.Range("D116").Value = "P31_Rel"
$ws.Range("D121").Value = "P32_Rel"
$ws.Range("D126").Value = "P33_Rel"
$ws.Range("D66").Value = "P20_Rel"
$ws.Range("D67").Value = "P20_Rel"
$ws.Range("D72").Value = "P21_Rel"
$ws.Range("D76").Value = "P22_Rel"
$ws.Range("D82").Value = "P23_Rel"
$ws.Range("D86").Value = "P24_Rel"
$ws.Range("D90").Value = "P25_Rel"
$ws.Range("D94").Value = "P26_Rel"

# Restore the sheet view scroll position and selection as recorded in the
# saved workbook state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 96
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J87").Select()
